$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.633.07"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.288.18"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "123.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0942"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("D16").Value = "2.634.30"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "2.290.08"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "43.700.85"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("E36").Value = "  +4.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.107"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E40").Value = "  +7.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "75.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "

# Row 47/48: swap ordi and TrustWalletToken with updated values
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +37.55%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.94%  "
